$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.84653948155325
$ws.Range("D2").Value = 8.910844356079908
$ws.Range("E2").Value = 15.92850303485443
$ws.Range("F2").Value = 39.01998439058193
$ws.Range("G2").Value = 46.50900457544447
$ws.Range("H2").Value = 18.46127169951713
$ws.Range("J2").Value = 11.6844617226592
$ws.Range("L2").Value = 12.05297736057702
$ws.Range("M2").Value = 18.47201947769856

$ws.Range("B3").Value = 19.53774408452444
$ws.Range("D3").Value = 8.815742727046159
$ws.Range("E3").Value = 15.7197458694059
$ws.Range("F3").Value = 39.18117701767221
$ws.Range("G3").Value = 46.46131388297334
$ws.Range("H3").Value = 18.52544070691826
$ws.Range("J3").Value = 11.60974676381659
$ws.Range("L3").Value = 11.80223657257622
$ws.Range("M3").Value = 18.25308771262805

$ws.Range("B4").Value = 19.34838208757212
$ws.Range("D4").Value = 8.756255722171353
$ws.Range("E4").Value = 15.58991593133557
$ws.Range("F4").Value = 39.29422538926107
$ws.Range("G4").Value = 46.45405802289123
$ws.Range("H4").Value = 18.57026407580117
$ws.Range("J4").Value = 11.5639513373112
$ws.Range("L4").Value = 11.64639935138364
$ws.Range("M4").Value = 18.11927595430853

$ws.Range("B5").Value = 19.27136255627308
$ws.Range("D5").Value = 8.731750361627721
$ws.Range("E5").Value = 15.53663287888061
$ws.Range("F5").Value = 39.34380778399441
$ws.Range("G5").Value = 46.45661706663344
$ws.Range("H5").Value = 18.58988647869719
$ws.Range("J5").Value = 11.54531904094648
$ws.Range("L5").Value = 11.58249678112594
$ws.Range("M5").Value = 18.06495088474836

$ws.Range("B6").Value = 19.25858494792694
$ws.Range("D6").Value = 8.727665655681417
$ws.Range("E6").Value = 15.52776369759192
$ws.Range("F6").Value = 39.3522524065922
$ws.Range("G6").Value = 46.4573743218571
$ws.Range("H6").Value = 18.59322647809321
$ws.Range("J6").Value = 11.54222723617338
$ws.Range("L6").Value = 11.5718640497228
$ws.Range("M6").Value = 18.05594404695713

$ws.Range("B7").Value = 19.3473426649883
$ws.Range("D7").Value = 8.755926287947956
$ws.Range("E7").Value = 15.58919880883804
$ws.Range("F7").Value = 39.29487987963222
$ws.Range("G7").Value = 46.45407023761566
$ws.Range("H7").Value = 18.57052322774323
$ws.Range("J7").Value = 11.56369992269176
$ws.Range("L7").Value = 11.64553904731974
$ws.Range("M7").Value = 18.11854241424372

$ws.Range("B8").Value = 19.74006892082435
$ws.Range("D8").Value = 8.878287101525331
$ws.Range("E8").Value = 15.85688975533542
$ws.Range("F8").Value = 39.07262737704752
$ws.Range("G8").Value = 46.48797782772143
$ws.Range("H8").Value = 18.48226728408884
$ws.Range("J8").Value = 11.65868923026077
$ws.Range("L8").Value = 11.96695679944675
$ws.Range("M8").Value = 18.39643464326359

$ws.Range("B9").Value = 20.50834290138347
$ws.Range("D9").Value = 9.10907264295537
$ws.Range("E9").Value = 16.36704215909221
$ws.Range("F9").Value = 38.74955949679745
$ws.Range("G9").Value = 46.72988691563872
$ws.Range("H9").Value = 18.35254229632689
$ws.Range("J9").Value = 11.84517266795913
$ws.Range("L9").Value = 12.57916592959819
$ws.Range("M9").Value = 18.94404390958547

$ws.Range("B10").Value = 21.06673212803965
$ws.Range("D10").Value = 9.272418272820323
$ws.Range("E10").Value = 16.73065471749251
$ws.Range("F10").Value = 38.58237555102235
$ws.Range("G10").Value = 47.01489267881765
$ws.Range("H10").Value = 18.28406953497409
$ws.Range("J10").Value = 11.98174210329314
$ws.Range("L10").Value = 13.01393286294966
$ws.Range("M10").Value = 19.34506522117598

$ws.Range("B11").Value = 21.31844136765538
$ws.Range("D11").Value = 9.345251498545514
$ws.Range("E11").Value = 16.89319932431407
$ws.Range("F11").Value = 38.52184159579334
$ws.Range("G11").Value = 47.16774325333802
$ws.Range("H11").Value = 18.2588290149945
$ws.Range("J11").Value = 12.04365434347299
$ws.Range("L11").Value = 13.20768635744595
$ws.Range("M11").Value = 19.52660552482877

$ws.Range("B12").Value = 21.41334594372605
$ws.Range("D12").Value = 9.372609168599281
$ws.Range("E12").Value = 16.95430343351328
$ws.Range("F12").Value = 38.50117216726137
$ws.Range("G12").Value = 47.22893869029481
$ws.Range("H12").Value = 18.25012689666875
$ws.Range("J12").Value = 12.06705784568167
$ws.Range("L12").Value = 13.28041808873253
$ws.Range("M12").Value = 19.59517270344213

$ws.Range("B13").Value = 21.39292613095069
$ws.Range("D13").Value = 9.366727277689478
$ws.Range("E13").Value = 16.94116405506581
$ws.Range("F13").Value = 38.50552314027752
$ws.Range("G13").Value = 47.21561219024106
$ws.Range("H13").Value = 18.25196288914523
$ws.Range("J13").Value = 12.06201948804482
$ws.Range("L13").Value = 13.26478331606923
$ws.Range("M13").Value = 19.580414254259

$ws.Range("B14").Value = 21.32625793727349
$ws.Range("D14").Value = 9.347506754318037
$ws.Range("E14").Value = 16.89823557653922
$ws.Range("F14").Value = 38.52009582590849
$ws.Range("G14").Value = 47.1727115704236
$ws.Range("H14").Value = 18.25809589754961
$ws.Range("J14").Value = 12.04558065209969
$ws.Range("L14").Value = 13.2136831656056
$ws.Range("M14").Value = 19.53225044528131

$ws.Range("B15").Value = 21.28536571533957
$ws.Range("D15").Value = 9.33570431266363
$ws.Range("E15").Value = 16.87188125335787
$ws.Range("F15").Value = 38.52931611189442
$ws.Range("G15").Value = 47.14686452693611
$ws.Range("H15").Value = 18.26196419028248
$ws.Range("J15").Value = 12.03550569681165
$ws.Range("L15").Value = 13.1822979891305
$ws.Range("M15").Value = 19.5027240310307

$ws.Range("B16").Value = 21.05022897372486
$ws.Range("D16").Value = 9.267627868617955
$ws.Range("E16").Value = 16.7199712863962
$ws.Range("F16").Value = 38.5866457651016
$ws.Range("G16").Value = 47.00536869112765
$ws.Range("H16").Value = 18.28583849847717
$ws.Range("J16").Value = 11.97769078482249
$ws.Range("L16").Value = 13.00118443051191
$ws.Range("M16").Value = 19.33317893452751

$ws.Range("B17").Value = 20.90533133045764
$ws.Range("D17").Value = 9.225479985761227
$ws.Range("E17").Value = 16.62601964515119
$ws.Range("F17").Value = 38.62580636292603
$ws.Range("G17").Value = 46.92449465756808
$ws.Range("H17").Value = 18.3020024411266
$ws.Range("J17").Value = 11.9421608982432
$ws.Range("L17").Value = 12.8890025608949
$ws.Range("M17").Value = 19.22890501432757

$ws.Range("B18").Value = 20.82177773137446
$ws.Range("D18").Value = 9.201099570179732
$ws.Range("E18").Value = 16.57171359206687
$ws.Range("F18").Value = 38.64978998094939
$ws.Range("G18").Value = 46.88016358056598
$ws.Range("H18").Value = 18.31185525289462
$ws.Range("J18").Value = 11.92170554545458
$ws.Range("L18").Value = 12.82410306166925
$ws.Range("M18").Value = 19.16884895442039

$ws.Range("B19").Value = 20.79345397118414
$ws.Range("D19").Value = 9.192821401386519
$ws.Range("E19").Value = 16.55328166514464
$ws.Range("F19").Value = 38.65816040140397
$ws.Range("G19").Value = 46.86552974034037
$ws.Range("H19").Value = 18.31528649524243
$ws.Range("J19").Value = 11.914776687686
$ws.Range("L19").Value = 12.80206661239624
$ws.Range("M19").Value = 19.14850275184768

$ws.Range("B20").Value = 20.9207785009713
$ws.Range("D20").Value = 9.229981068623349
$ws.Range("E20").Value = 16.63604891437178
$ws.Range("F20").Value = 38.62148643759664
$ws.Range("G20").Value = 46.93287777563841
$ws.Range("H20").Value = 18.30022419784119
$ws.Range("J20").Value = 11.94594521036121
$ws.Range("L20").Value = 12.90098381713875
$ws.Range("M20").Value = 19.24001384989732

$ws.Range("B21").Value = 21.34585181750485
$ws.Range("D21").Value = 9.353158416216381
$ws.Range("E21").Value = 16.91085715063471
$ws.Range("F21").Value = 38.51575414845956
$ws.Range("G21").Value = 47.18522278915757
$ws.Range("H21").Value = 18.25627120615321
$ws.Range("J21").Value = 12.05041034201367
$ws.Range("L21").Value = 13.22871030399646
$ws.Range("M21").Value = 19.54640255770231

$ws.Range("B22").Value = 21.62122384220309
$ws.Range("D22").Value = 9.43235913767472
$ws.Range("E22").Value = 17.08783316855376
$ws.Range("F22").Value = 38.45979614031553
$ws.Range("G22").Value = 47.36944699958602
$ws.Range("H22").Value = 18.23253701870721
$ws.Range("J22").Value = 12.11843864733419
$ws.Range("L22").Value = 13.43915192114229
$ws.Range("M22").Value = 19.7455840400023

$ws.Range("B23").Value = 21.47450149890969
$ws.Range("D23").Value = 9.390210911776222
$ws.Range("E23").Value = 16.99362966669968
$ws.Range("F23").Value = 38.4884523987688
$ws.Range("G23").Value = 47.26936622126458
$ws.Range("H23").Value = 18.24474566477574
$ws.Range("J23").Value = 12.0821566161843
$ws.Range("L23").Value = 13.32719678762981
$ws.Range("M23").Value = 19.63939066396748

$ws.Range("B24").Value = 20.91379560775219
$ws.Range("D24").Value = 9.227946592453575
$ws.Range("E24").Value = 16.63151558852948
$ws.Range("F24").Value = 38.62343489982077
$ws.Range("G24").Value = 46.92908102451601
$ws.Range("H24").Value = 18.30102639800928
$ws.Range("J24").Value = 11.94423441159398
$ws.Range("L24").Value = 12.89556834699462
$ws.Range("M24").Value = 19.2349918766058

$ws.Range("B25").Value = 20.30122030251426
$ws.Range("D25").Value = 9.047685575239786
$ws.Range("E25").Value = 16.23084476237261
$ws.Range("F25").Value = 38.82473405697056
$ws.Range("G25").Value = 46.64559765273896
$ws.Range("H25").Value = 18.38295231258419
$ws.Range("J25").Value = 11.7947599455838
$ws.Range("L25").Value = 12.41591045039407
$ws.Range("M25").Value = 18.79592034185512

